$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '71.791.93'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '4.031.11'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.00%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.56'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.59'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.25%  '
$ws.Range('E7').Value = '  +1.25%  '
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.744'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.177'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000345'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.91'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +9.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.95'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.669.17'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.080.95'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.46'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.38'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.30%  '
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('E19').Value = '  -1.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.755.57'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '447.04'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.63'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '94.93'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.45'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.20'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.07'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.12'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.37'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.58%  '
$ws.Range('E29').Value = '  +3.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '701.75'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  +3.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.07'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +17.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.91'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '68.11'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0910'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.450'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '41.66'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.12%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.154'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.25%  '
$ws.Range('B39').Value = 'ThetaToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.53'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +17.69%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +2.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.55'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.40%  '
$ws.Range('E46').Value = '  +3.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.18'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.74%  '
$ws.Range('E48').Value = '  +7.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000280'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +18.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.39'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0344'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.47%  '
